$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: "Change J5-LCD-3V3" -> "Change J5-LCD-3V3 to DNP"
$ws.Range("B5").Value = "Change J5-LCD-3V3 to DNP"

# Update selection to B6 on Sheet1 (as recorded in the saved file)
$ws.Activate()
$ws.Range("B6").Select()
